$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update datetimes in column A (rows 3-11)
$ws.Range("A3").Value = "2022-09-14 22:23:42"
$ws.Range("A4").Value = "2022-09-14 22:23:43"
$ws.Range("A5").Value = "2022-09-14 22:23:45"
$ws.Range("A6").Value = "2022-09-14 22:23:53"
$ws.Range("A7").Value = "2022-09-14 22:24:04"
$ws.Range("A8").Value = "2022-09-14 22:24:18"
$ws.Range("A9").Value = "2022-09-14 22:24:53"
$ws.Range("A10").Value = "2022-09-14 22:24:56"
$ws.Range("A11").Value = "2022-09-14 22:25:12"

# Update US name in column B (rows 3-11)
$ws.Range("B3").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B4").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B5").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B6").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B7").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B8").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B9").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B10").Value = "MER_CT_ChamanculoCS_37"
$ws.Range("B11").Value = "MER_CT_ChamanculoCS_37"

# Remove row 12 entirely (shifts remaining rows up, none below it here)
$ws.Rows.Item(12).Delete()
